$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18002.533
$ws.Range("I21").Value = 25004.75
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 25004.75
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -24536.75
$ws.Range("N21").Value = -10936
$ws.Range("H23").Value = 18002.533
$ws.Range("I23").Value = 25004.75
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 25004.75
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -24770.75
$ws.Range("N23").Value = -10468
$ws.Range("H125").Value = 1449.7778
$ws.Range("I125").Value = 1519
$ws.Range("J125").Value = 1394.4
$ws.Range("K125").Value = 13671
$ws.Range("L125").Value = 12549.6
$ws.Range("M125").Value = -11211
$ws.Range("N125").Value = -17469.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("H140").Value = 32966.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 32966.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 32966.668
$ws.Range("N140").Value = -43326.668
$ws.Range("H141").Value = 57500
$ws.Range("J141").Value = 57500
$ws.Range("L141").Value = 57500
$ws.Range("N141").Value = -67860
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("M140").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 14400.5
$ws.Range("I97").Value = 5925.2
$ws.Range("K97").Value = 5925.2
$ws.Range("M97").Value = -4934.2
$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -23134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1070.83
$ws.Range("I31").Value = 696.6667
$ws.Range("J31").Value = 1566.814
$ws.Range("K31").Value = 696.6667
$ws.Range("L31").Value = 1566.814
$ws.Range("M31").Value = -401.6667
$ws.Range("N31").Value = -2156.814
$ws.Range("H34").Value = 1070.83
$ws.Range("I34").Value = 696.6667
$ws.Range("J34").Value = 1566.814
$ws.Range("K34").Value = 696.6667
$ws.Range("L34").Value = 1566.814
$ws.Range("M34").Value = -494.6667
$ws.Range("N34").Value = -1970.814
$ws.Range("H140").Value = 36753.453
$ws.Range("J140").Value = 36753.453
$ws.Range("L140").Value = 36753.453
$ws.Range("N140").Value = -47113.453
$ws.Range("H141").Value = 7678.2
$ws.Range("J141").Value = 7678.2
$ws.Range("L141").Value = 7678.2
$ws.Range("N141").Value = -18038.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1503.08
$ws.Range("I114").Value = 1160.2858
$ws.Range("J114").Value = 1939.3636
$ws.Range("K114").Value = 3480.8574
$ws.Range("L114").Value = 5818.0908
$ws.Range("M114").Value = -226.8574000000003
$ws.Range("N114").Value = -12326.0908

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5128.4136
$ws.Range("I80").Value = 3990.4211
$ws.Range("J80").Value = 7290.6
$ws.Range("K80").Value = 3990.4211
$ws.Range("L80").Value = 7290.6
$ws.Range("M80").Value = -2992.4211
$ws.Range("N80").Value = -9286.6
$ws.Range("H83").Value = 5128.4136
$ws.Range("I83").Value = 3990.4211
$ws.Range("J83").Value = 7290.6
$ws.Range("K83").Value = 19952.1055
$ws.Range("L83").Value = 36453
$ws.Range("M83").Value = -14960.1055
$ws.Range("N83").Value = -46437
$ws.Range("H116").Value = 44842
$ws.Range("J116").Value = 44842
$ws.Range("L116").Value = 44842
$ws.Range("N116").Value = -54020

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9111.111000000001
$ws.Range("H7").Value = 71431840
$ws.Range("I7").Value = 125001790
$ws.Range("J7").Value = 5235
$ws.Range("K7").Value = 125001790
$ws.Range("L7").Value = 5235
$ws.Range("M7").Value = -125001678
$ws.Range("N7").Value = -5459
$ws.Range("H126").Value = 71431840
$ws.Range("I126").Value = 125001790
$ws.Range("J126").Value = 5235
$ws.Range("K126").Value = 375005370
$ws.Range("L126").Value = 15705
$ws.Range("M126").Value = -375002900
$ws.Range("N126").Value = -20645
$ws.Range("H132").Value = 3512.879
$ws.Range("I132").Value = 3027.9048
$ws.Range("J132").Value = 4361.5835
$ws.Range("K132").Value = 9083.714399999999
$ws.Range("L132").Value = 13084.7505
$ws.Range("M132").Value = -6553.714399999999
$ws.Range("N132").Value = -18144.7505
$ws.Range("H135").Value = 30411.875
$ws.Range("J135").Value = 30411.875
$ws.Range("L135").Value = 30411.875
$ws.Range("N135").Value = -40551.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 70005
$ws.Range("I24").Value = 70000
$ws.Range("J24").Value = 70010
$ws.Range("K24").Value = 70000
$ws.Range("L24").Value = 70010
$ws.Range("M24").Value = -69770
$ws.Range("N24").Value = -70470
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H81").Value = 11906895
$ws.Range("I81").Value = 14287674
$ws.Range("K81").Value = 28575348
$ws.Range("M81").Value = -28574287
$ws.Range("H84").Value = 11906895
$ws.Range("I84").Value = 14287674
$ws.Range("K84").Value = 142876740
$ws.Range("M84").Value = -142871436
$ws.Range("H107").Value = 9092045
$ws.Range("I107").Value = 1106.3334
$ws.Range("J107").Value = 20001170
$ws.Range("K107").Value = 3319.0002
$ws.Range("L107").Value = 60003510
$ws.Range("M107").Value = -1399.0002
$ws.Range("N107").Value = -60007350
$ws.Range("H136").Value = 401227.88
$ws.Range("I136").Value = 527405.2
$ws.Range("K136").Value = 1582215.6
$ws.Range("M136").Value = -1579665.6
$ws.Range("H137").Value = 64618.832
$ws.Range("J137").Value = 64618.832
$ws.Range("L137").Value = 64618.832
$ws.Range("N137").Value = -74818.83199999999
$ws.Range("H138").Value = 46000
$ws.Range("J138").Value = 46000
$ws.Range("L138").Value = 46000
$ws.Range("N138").Value = -56280
$ws.Range("H139").Value = 49621.875
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49621.875
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49621.875
$ws.Range("N139").Value = -59901.875
$ws.Range("H140").Value = 40071.46
$ws.Range("J140").Value = 40071.46
$ws.Range("L140").Value = 40071.46
$ws.Range("N140").Value = -50431.46
$ws.Range("H141").Value = 27878.143
$ws.Range("J141").Value = 27878.143
$ws.Range("L141").Value = 27878.143
$ws.Range("N141").Value = -38238.143
$ws.Range("N25").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("N30").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("M139").ClearContents()
